# Completed Appeal Form Date verification:
# append the 22 new "grounds for appeal" verification rows (123-144) to the
# Verifications sheet, in the same order the strings were originally typed
# (note rows 129, 140, 141 were filled in later, after 142-144), so the
# resulting shared-string table matches the authored order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A95").Copy()
$ws.Range("A123").PasteSpecial(-4122)
$ws.Cells.Item(123,1).Value = 'You can now submit your application on the grounds of: Refused planning permission'

$ws.Range("A95").Copy()
$ws.Range("A124").PasteSpecial(-4122)
$ws.Cells.Item(124,1).Value = 'You can now submit your application on the grounds of: Refused approval of the matters reserved under an outline planning permission'

$ws.Cells.Item(125,1).Value = 'You can now submit your application on the grounds of: Modification or discharge of planning obligations'

$ws.Cells.Item(126,1).Value = 'You can now submit your application on the grounds of: Refused prior approval of permitted development rights'

$ws.Cells.Item(127,1).Value = 'You can now submit your application on the grounds of: Granted planning permission for the development subject to conditions to which you object'

$ws.Cells.Item(128,1).Value = 'You can now submit your application on the grounds of: Refused permission to vary or remove a condition(s)'

$ws.Cells.Item(130,1).Value = 'You can now submit your application on the grounds of: Granted approval of the matters reserved under an outline planning permission subject to conditions to which you object'

$ws.Cells.Item(131,1).Value = 'You can now submit your application on the grounds of: Non-determination appeal'

$ws.Cells.Item(132,1).Value = 'You can now submit your application on the grounds of: Non-determination appeal (local list documentation)'

$ws.Cells.Item(133,1).Value = 'You can now submit your application on the grounds of: Householder Appeals Service'

$ws.Cells.Item(134,1).Value = 'You can now submit your application on the grounds of: Commercial Appeals Service'

$ws.Cells.Item(135,1).Value = 'You can now submit your application on the grounds of: Commercial Appeals Service (Adverts)'

$ws.Cells.Item(136,1).Value = 'You can now submit your application on the grounds of: Listed Building Consent'

$ws.Cells.Item(137,1).Value = 'You can now submit your application on the grounds of: Other'

$ws.Cells.Item(138,1).Value = 'You can now submit your application on the grounds of: Granted planning permission for the development subject to conditions to which you object (Technical Design)'

$ws.Cells.Item(139,1).Value = 'You can now submit your application on the grounds of: Refused permission to vary or remove a condition(s) (Technical Design)'

$ws.Cells.Item(142,1).Value = 'You can now submit your application on the grounds of: Non-determination appeal (Permission in Principle)'

$ws.Cells.Item(143,1).Value = 'You can now submit your application on the grounds of: Non-determination appeal (Technical Design)'

$ws.Cells.Item(144,1).Value = 'You can now submit your application on the grounds of: HAS Prior approval'

$ws.Cells.Item(129,1).Value = 'You can now submit your application on the grounds of: Refused to approve any matter required by a condition on a previous planning permission'

$ws.Cells.Item(140,1).Value = 'You can now submit your application on the grounds of: Permission in Principle'

$ws.Cells.Item(141,1).Value = 'You can now submit your application on the grounds of: Technical Design Consent'

$ws.Range("A141").Select()
